$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORTDATE
$ws.Range("H2").Value = "2019-12-31 00:00:00"

# BASIC_EPS, DEDUCT_BASIC_EPS
$ws.Range("I2").Value = 0.73
$ws.Range("J2").Value = 0.71

# TOTAL_OPERATE_INCOME, PARENT_NETPROFIT
$ws.Range("K2").Value = 779287856.16
$ws.Range("L2").Value = 262013571.95

# WEIGHTAVG_ROE
$ws.Range("M2").Value = 36.55

# YSTZ, SJLTZ (previously blank, now populated numbers)
$ws.Range("N2").Value = 6.2308604641
$ws.Range("O2").Value = 14.7278094268

# BPS, MGJYXJJE, XSMLL
$ws.Range("P2").Value = 8.336805892556001
$ws.Range("Q2").Value = 2.913204943444
$ws.Range("R2").Value = 59.4973922916

# ISNEW (text "0", not number)
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "0"

# QDATE, DATATYPE
$ws.Range("AC2").Value = "2019Q4"
$ws.Range("AD2").Value = "2019年 年报"

# DATAYEAR (text "2019", not number)
$ws.Range("AE2").NumberFormat = "@"
$ws.Range("AE2").Value = "2019"

# DATEMMDD
$ws.Range("AF2").Value = "年报"
